# Daily attendance processing - 2026-01-09 05:14:51
# Swap the order of "System" and the recorded-by email address in the
# "Recorded By" column (G) wherever both appear together.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$oldVal = "System, dnasr281@gmail.com"
$newVal = "dnasr281@gmail.com, System"

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    if ($cell.Text -eq $oldVal) {
        $cell.Value = $newVal
    }
}
